## Adds the "thesis / ARPA errors difference" table (Table2) to Sheet1.
## Table2 lives at A9:D13 with columns: transcription, recording_id,
## model_transcription, difference - four rows of ASR error-analysis data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 9)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "transcription"
$ws.Range("B9").Value = "recording_id"
$ws.Range("C9").Value = "model_transcription"
$ws.Range("D9").Value = "difference"

# ---------------------------------------------------------------------
# 2. Data rows (10-13)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "golf bravo victor juliett india is identified good afternoon"
$ws.Range("B10").Value = "051_0227"
$ws.Range("C10").Value = "goll bravo victor juliett india is identified good afternoon"
$ws.Range("D10").Value = "{'golf'}"

$ws.Range("A11").Value = "exact"
$ws.Range("B11").Value = "111_0627"
$ws.Range("C11").Value = "expect"
$ws.Range("D11").Value = "{'exact'}"

$ws.Range("A12").Value = "roger what is your position"
$ws.Range("B12").Value = "101_0199"
$ws.Range("C12").Value = "roger ah what is your position"
$ws.Range("D12").Value = "set()"

$ws.Range("A13").Value = "japan air four one nine contact milan one three four five two bye"
$ws.Range("B13").Value = "101_0308"
$ws.Range("C13").Value = "german air four one nine contact milan one three four five two bye"
$ws.Range("D13").Value = "{'japan'}"

# ---------------------------------------------------------------------
# 3. Turn A9:D13 into a native Excel table ("Table2")
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A9:D13"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table2"

# ---------------------------------------------------------------------
# 4. Column width tweaks (directory / filename columns got wider)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.6
$ws.Columns.Item(3).ColumnWidth = 20.3

# ---------------------------------------------------------------------
# 5. Make the new table's range the active selection, like the author
#    left it selected after inserting it.
# ---------------------------------------------------------------------
$ws.Range("A9:D13").Select()
